$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AUR-UM192"
$ws.Range("A4").Value = "AUR-UM192"

$ws.Range("A4").Select() | Out-Null
